# Edit script: CUZ Quarterly Financials update
# - Inserts two new columns (new D & E) before existing data, shifting old D:K -> F:M
# - Populates the new D/E columns with the latest two quarter-end figures
# - Applies a handful of corrected historical figures noted in the source update

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at D:E; Excel shifts the old D:K block to F:M
$ws.Columns("D:E").Insert()

# Copy number formatting from the (now-shifted) first data column into the new D:E columns
# Period-ending date header rows use the date style; all other rows use the numeric style
$ws.Range("F7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("F38").Copy()
$ws.Range("D38:E38").PasteSpecial(-4122)
$ws.Range("F80").Copy()
$ws.Range("D80:E80").PasteSpecial(-4122)

$ws.Range("F8:F35").Copy()
$ws.Range("D8:E35").PasteSpecial(-4122)
$ws.Range("F39:F77").Copy()
$ws.Range("D39:E77").PasteSpecial(-4122)
$ws.Range("F81:F102").Copy()
$ws.Range("D81:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new D (latest quarter) and E (prior quarter) columns
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 122700
$ws.Range("E8").Value = 118700
$ws.Range("D9").Value = 42300
$ws.Range("E9").Value = 41700
$ws.Range("D10").Value = 80400
$ws.Range("E10").Value = 77000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = -100
$ws.Range("D15").Value = 45500
$ws.Range("E15").Value = 45100
$ws.Range("D17").Value = 92100
$ws.Range("E17").Value = 91500
$ws.Range("D18").Value = 30600
$ws.Range("E18").Value = 27200
$ws.Range("D20").Value = 500
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 76600
$ws.Range("E21").Value = 72200
$ws.Range("D22").Value = 10400
$ws.Range("E22").Value = 9600
$ws.Range("D23").Value = 20700
$ws.Range("E23").Value = 17600
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 20700
$ws.Range("E26").Value = 17600
$ws.Range("D27").Value = 22400
$ws.Range("E27").Value = 19500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -500
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 22400
$ws.Range("E33").Value = 19500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 22400
$ws.Range("E35").Value = 19500
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 2500
$ws.Range("E41").Value = 82700
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 96900
$ws.Range("E43").Value = 86900
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 5100
$ws.Range("E45").Value = 5000
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 161900
$ws.Range("E47").Value = 154100
$ws.Range("D48").Value = 3714700
$ws.Range("E48").Value = 3645700
$ws.Range("D49").Value = 145900
$ws.Range("E49").Value = 155000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 13400
$ws.Range("E52").Value = 8000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 4146300
$ws.Range("E54").Value = 4143600
$ws.Range("D57").Value = 110200
$ws.Range("E57").Value = 114200
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 1062600
$ws.Range("E61").Value = 1065000
$ws.Range("D62").Value = "NA"
$ws.Range("E62").Value = "NA"
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1380400
$ws.Range("E66").Value = 1373400
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 6900
$ws.Range("E70").Value = 6900
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -1129400
$ws.Range("E72").Value = -1124500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 2759000
$ws.Range("E76").Value = 2763400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 22400
$ws.Range("E81").Value = 19500
$ws.Range("D83").Value = 45500
$ws.Range("E83").Value = 45100
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 51000
$ws.Range("E89").Value = 73300
$ws.Range("D91").Value = -91200
$ws.Range("E91").Value = -32700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -102100
$ws.Range("E94").Value = -48000
$ws.Range("D96").Value = -27300
$ws.Range("E96").Value = -27300
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -29300
$ws.Range("E100").Value = -52700
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -80400
$ws.Range("E102").Value = -27400

# Corrected figures in the historical columns (not a pure shift of the old data)
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("I91").Value = -78700
$ws.Range("J91").Value = -80400
$ws.Range("H94").Value = 119500
$ws.Range("I94").Value = -81400

